$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new cells
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin box border, centered horizontally, top vertically
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2

# Apply the exact same style to A2 via copy/paste-format so that both cells
# end up referencing the same cellXf entry (avoids creating a duplicate style).
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)
